$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "client"/"hours" columns (A,B) two columns to the right
# (C,D), freeing up A,B for the new indicator columns. This also carries the
# custom column width that lived on column A over to column C, matching it
# exactly since the underlying column-properties move with the column.
$null = $ws.Range("A1:B1").EntireColumn.Insert()

# New header row: icow, ihorse, iclient, ihours, itonto
$ws.Range("A1").Value = "icow"
$ws.Range("B1").Value = "ihorse"
$ws.Range("C1").Value = "iclient"
$ws.Range("D1").Value = "ihours"
$ws.Range("E1").Value = "itonto"

# Row 2 (was American/200)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("E2").Value = 1

# Row 3 (was BA/100)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 1

# Row 4 (was AirFrance/500)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("E4").Value = 1

# Move the selection to C2, matching the new active cell
$null = $ws.Range("C2").Select()
